$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.934.60"
$ws.Range("E2").Value = "  -3.20%  "

$ws.Range("D3").Value = "3.074.17"
$ws.Range("E3").Value = "  -2.08%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "544.92"
$ws.Range("E5").Value = "  -3.87%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.01"
$ws.Range("E6").Value = "  -8.63%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("D8").Value = "3.069.74"
$ws.Range("E8").Value = "  -1.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.489"
$ws.Range("E9").Value = "  -1.75%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.49"
$ws.Range("E10").Value = "  -5.99%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  -0.82%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  -1.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "34.76"
$ws.Range("E13").Value = "  -2.89%  "

$ws.Range("E14").Value = "  -2.81%  "

$ws.Range("D15").Value = "3.569.37"
$ws.Range("E15").Value = "  -2.23%  "

$ws.Range("D16").Value = "63.021.83"
$ws.Range("E16").Value = "  -3.19%  "

$ws.Range("E17").Value = "  -1.37%  "

$ws.Range("D18").Value = "3.073.19"
$ws.Range("E18").Value = "  -2.26%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.60"
$ws.Range("E19").Value = "  -1.81%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "482.00"
$ws.Range("E20").Value = "  -7.98%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.39"
$ws.Range("E21").Value = "  -3.16%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.696"
$ws.Range("E22").Value = "  -0.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.08"
$ws.Range("E23").Value = "  -4.51%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "77.45"
$ws.Range("E24").Value = "  -1.30%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.17"
$ws.Range("E25").Value = "  -4.36%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.71"
$ws.Range("E27").Value = "  -3.42%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  -5.38%  "

$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.48"
$ws.Range("E34").Value = "  -6.36%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "520.59"
$ws.Range("E35").Value = "  -7.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -1.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.11"
$ws.Range("E37").Value = "  -5.54%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0398"
$ws.Range("E38").Value = "  -9.74%  "

$ws.Range("D39").Value = "3.078.40"
$ws.Range("E39").Value = "  +0.21%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0786"
$ws.Range("E40").Value = "  -3.38%  "

$ws.Range("E41").Value = "  -2.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.04"
$ws.Range("E42").Value = "  -2.31%  "

$ws.Range("E43").Value = "  -7.86%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.252"
$ws.Range("E44").Value = "  -1.67%  "

$ws.Range("E45").Value = "  +0.01%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.04"
$ws.Range("E46").Value = "  -6.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "120.51"
$ws.Range("E47").Value = "  +1.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "24.17"
$ws.Range("E48").Value = "  -2.93%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.107"
$ws.Range("E49").Value = "  -1.80%  "

$ws.Range("D50").Value = "0.0₃0498"
$ws.Range("E50").Value = "  -4.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.32"
$ws.Range("E51").Value = "  +60.20%  "

$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.93"
$ws.Range("E30").Value = "  -9.07%  "

$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.23"
$ws.Range("E31").Value = "  -0.03%  "

$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.13"
$ws.Range("E32").Value = "  -1.60%  "

$ws.Range("B33").Value = "OKB"
$ws.Range("C33").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "61.16"
$ws.Range("E33").Value = "  +15.66%  "

